$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# -----------------------------------------------------------------
# 1) Trim the tail of the "Momentum Storage" paragraph: keep up to
#    "...How much depends on orbital altitude. " and drop the rest
#    (the old text about the 200km circular orbit / aerobraking /
#    desaturations), replacing it with a brand-new run so the
#    preceding " disturbances" run is left untouched.
# -----------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" this is not cyclic and will accumulate throughout the orbit. How much depends on orbital altitude. The final circular 200km circular orbit will drive the design, however it is noted that immediately after aerobraking the orbit will be lower and will require more frequent momentum desaturations. Frequency of desaturations is TBD.")
$momentumPara = $rng.Paragraphs(1)
$rng.Delete()
$rng.InsertAfter(" this is not cyclic and will accumulate throughout the orbit. How much depends on orbital altitude. ")

# -----------------------------------------------------------------
# 2) Insert the new paragraphs right after the "Momentum Storage"
#    paragraph, one by one, carrying the exact pPr/run content.
# -----------------------------------------------------------------
$cur = $momentumPara

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:spacing w:after="0"/><w:ind w:left="720"/></w:pPr><w:r><w:t>The worst case scenario is considered, taking into account the torques needed when driving the solar arrays and HGA</w:t></w:r></w:p>')

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>')

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Appendages:</w:t></w:r></w:p>')

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Solar array:</w:t></w:r></w:p>')

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Circular shape, rotation only about radial axis</w:t></w:r></w:p>')

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>HGA:</w:t></w:r></w:p>')

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Cylindrical shape. </w:t></w:r></w:p>')

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Gimbal located on edge of body</w:t></w:r></w:p>')

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Gimbal mechanism 45kg, see designParams.xlsx first sheet [Deep Space </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Comms</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 2016, p210]</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# -----------------------------------------------------------------
# 3) Remove the old "_GoBack" bookmark that used to sit on the empty
#    bullet paragraph after "Landing radar (Terminal Descent Sensor
#    TDS)" -- it has now moved to the paragraph above.
# -----------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Landing radar (Terminal Descent Sensor TDS)")
$landingPara = $rng2.Paragraphs(1)
$bookmarkPara = $landingPara.Next()
$bookmarkPara.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr></w:p>')

# -----------------------------------------------------------------
# 4) Move the "lastRenderedPageBreak" marker from the "Budgets"
#    heading to the "CDH:" heading.
# -----------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Budgets")
$budgetsPara = $rng3.Paragraphs(1)
$budgetsPara.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Heading2"/><w:spacing w:before="0"/></w:pPr><w:r><w:t>Budgets</w:t></w:r><w:r><w:t xml:space="preserve"> (Basic sizing </w:t></w:r><w:r><w:t xml:space="preserve">based on MRO and MSL, </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>NO CALCULATIONS</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>')

$rng4 = $d.Content
$rng4.Find.Execute("CDH:")
$cdhPara = $rng4.Paragraphs(1)
$cdhPara.Range.InsertXML('<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Heading4"/><w:spacing w:before="0"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>CDH:</w:t></w:r></w:p>')

Write-Host "edit complete"
